# Swap the data (columns B through AD) between two adjacent rows,
# leaving column A (the running index) untouched.
function Swap-RowData($ws, $Row1, $Row2) {
    $range1 = $ws.Range("B" + $Row1 + ":AD" + $Row1)
    $range2 = $ws.Range("B" + $Row2 + ":AD" + $Row2)

    $data1 = $range1.Value()
    $data2 = $range2.Value()

    $range1.Value() = $data2
    $range2.Value() = $data1
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Swap-RowData $ws 120 121
Swap-RowData $ws 151 152
Swap-RowData $ws 161 162
